$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memberlist")

for ($i = 1; $i -le 20; $i++) {
    $num = "{0:D3}" -f $i
    $ws.Cells.Item($i + 1, 2).Value = "upload/$num.JPG"
}

$ws.Activate()
$ws.Range("B9").Select()
